# Generate Report for Archive
#
# Rows 4 and 5 on every sheet (Overview, zh-cn, de-de) swap file identity:
# "fe2645a6-...md" is promoted to row 4 (ahead of "f7783899-...md", which
# moves down to row 5). The promoted row 4 also flips its handoff status
# from "Ready for handoff" to "In Translation"; row 5 keeps
# "Ready for handoff". The zh-cn/de-de sheets carry the matching handoff
# file name + handoff datetime along with each row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview": A = file name, B = zh-cn status, C = de-de status
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A4").Value = "fe2645a6-7783-4020-a687-f4d97400839a.md"
$ov.Range("B4").Value = "In Translation"
$ov.Range("C4").Value = "In Translation"

$ov.Range("A5").Value = "f7783899-1201-4321-83dd-2675893f7f26.md"
$ov.Range("B5").Value = "Ready for handoff"
$ov.Range("C5").Value = "Ready for handoff"

# ---------------------------------------------------------------
# Sheet "zh-cn": A = file, B = status, C = handoff file, D = handoff datetime
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A4").Value = "fe2645a6-7783-4020-a687-f4d97400839a.md"
$zh.Range("B4").Value = "In Translation"
$zh.Range("C4").Value = "fe2645a6-7783-4020-a687-f4d97400839a.b61eb865429a80a8d793f3a896ec316f60084507.zh-cn.xlf"
$zh.Range("D4").Value = "2016-02-24 08:35:13"

$zh.Range("A5").Value = "f7783899-1201-4321-83dd-2675893f7f26.md"
$zh.Range("B5").Value = "Ready for handoff"
$zh.Range("C5").Value = "f7783899-1201-4321-83dd-2675893f7f26.862c3495ceee30cdf24ea0815fe639aab36279db.zh-cn.xlf"
$zh.Range("D5").Value = "2016-02-24 08:35:54"

# ---------------------------------------------------------------
# Sheet "de-de": A = file, B = status, C = handoff file, D = handoff datetime
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A4").Value = "fe2645a6-7783-4020-a687-f4d97400839a.md"
$de.Range("B4").Value = "In Translation"
$de.Range("C4").Value = "fe2645a6-7783-4020-a687-f4d97400839a.b61eb865429a80a8d793f3a896ec316f60084507.de-de.xlf"
$de.Range("D4").Value = "2016-02-24 08:35:23"

$de.Range("A5").Value = "f7783899-1201-4321-83dd-2675893f7f26.md"
$de.Range("B5").Value = "Ready for handoff"
$de.Range("C5").Value = "f7783899-1201-4321-83dd-2675893f7f26.862c3495ceee30cdf24ea0815fe639aab36279db.de-de.xlf"
$de.Range("D5").Value = "2016-02-24 08:36:04"
